# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

$rushing   = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet updates ---
# Row 2 - A.Dalton
$rushing.Range("D2").Value = 4
$rushing.Range("E2").Value = 2
$rushing.Range("F2").Value = 3

# Row 5 - D.Montgomery
$rushing.Range("C5").Value = 150
$rushing.Range("D5").Value = 77
$rushing.Range("E5").Value = 30
$rushing.Range("F5").Value = 37

# Row 7 - K.Herbert
$rushing.Range("D7").Value = 10
$rushing.Range("E7").Value = 3

# Row 9 - A.Robinson
$rushing.Range("C9").Value = 2
$rushing.Range("F9").Value = 1

# --- Receiving sheet updates ---
# Row 2 - D.Montgomery
$receiving.Range("C2").Value = 57
$receiving.Range("D2").Value = 49

# Row 3 - D.Williams
$receiving.Range("C3").Value = 12

# Row 5 - A.Robinson
$receiving.Range("C5").Value = 46
$receiving.Range("D5").Value = 31
$receiving.Range("E5").Value = 17

# Row 6 - D.Mooney
$receiving.Range("C6").Value = 92
$receiving.Range("D6").Value = 55
$receiving.Range("E6").Value = 32
$receiving.Range("G6").Value = 10
$receiving.Range("H6").Value = 5

# Row 7 - M.Goodwin
$receiving.Range("C7").Value = 26
$receiving.Range("G7").Value = 3

# Row 8 - D.Byrd
$receiving.Range("C8").Value = 26
$receiving.Range("E8").Value = 6
$receiving.Range("F8").Value = 3

# Row 11 - C.Kmet
$receiving.Range("C11").Value = 77
$receiving.Range("D11").Value = 50
$receiving.Range("G11").Value = 13
$receiving.Range("H11").Value = 6

# Row 12 - J.Graham
$receiving.Range("E12").Value = 3

# Row 13 - J.James
$receiving.Range("C13").Value = 8
$receiving.Range("D13").Value = 7

$wb.Save()
